# Update workbook from v0.2 to v1.0.1 and fix TC2/TC3 step ordering.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# 1. Update the version number cell (D2) from "0.1" to "1.0.1"
$ws.Range("D2").Value = "1.0.1"

# 2. Swap the step/result content between the TC2 block (row 20) and the
#    TC3 block (row 28): the "realizar a liquidacao" step should move to
#    TC3, and the "atribuir/desatribuir" step should move to TC2.
$tc2Step = $ws.Range("B20").Value()
$tc2Result = $ws.Range("D20").Value()
$tc3Step = $ws.Range("B28").Value()
$tc3Result = $ws.Range("D28").Value()

$ws.Range("B20").Value = $tc3Step
$ws.Range("D20").Value = $tc3Result
$ws.Range("B28").Value = $tc2Step
$ws.Range("D28").Value = $tc2Result
